$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 400
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 400
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 400
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -538

$ws.Range("H70").Value = 3300
$ws.Range("J70").Value = 3233.3333
$ws.Range("L70").Value = 9699.999899999999
$ws.Range("N70").Value = -10239.9999

$ws.Range("H73").Value = 3300
$ws.Range("J73").Value = 3233.3333
$ws.Range("L73").Value = 9699.999899999999
$ws.Range("N73").Value = -11571.9999

$ws.Range("H81").Value = 35000
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -36996

$ws.Range("H84").Value = 35000
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -114984

$ws.Range("H107").Value = 1094.0952
$ws.Range("I107").Value = 1805.7142
$ws.Range("J107").Value = 738.2857
$ws.Range("K107").Value = 1805.7142
$ws.Range("L107").Value = 738.2857
$ws.Range("M107").Value = 114.2858000000001
$ws.Range("N107").Value = -4578.2857

$ws.Range("H112").Value = 2661.111
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2661.111
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 7983.333
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -10199.333

$ws.Range("H115").Value = 3409.6
$ws.Range("I115").Value = 3626.4614
$ws.Range("K115").Value = 10879.3842
$ws.Range("M115").Value = -9312.3842

$ws.Range("H118").Value = 650
$ws.Range("I118").Value = 650
$ws.Range("K118").Value = 1950
$ws.Range("M118").Value = -293

$ws.Range("H127").Value = 1013.64703
$ws.Range("I127").Value = 360
$ws.Range("J127").Value = 1947.4286
$ws.Range("K127").Value = 1080
$ws.Range("L127").Value = 5842.2858
$ws.Range("M127").Value = 3880
$ws.Range("N127").Value = -15762.2858

$ws.Range("H129").Value = 728.1429000000001
$ws.Range("I129").Value = 486.75
$ws.Range("J129").Value = 1050
$ws.Range("K129").Value = 1460.25
$ws.Range("L129").Value = 3150
$ws.Range("M129").Value = 3539.75
$ws.Range("N129").Value = -13150

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 920850.25
$ws.Range("I2").Value = 1566.55
$ws.Range("J2").Value = 2452989.8
$ws.Range("K2").Value = 1566.55
$ws.Range("L2").Value = 2452989.8
$ws.Range("M2").Value = -1453.55
$ws.Range("N2").Value = -2453215.8

$ws.Range("H32").Value = 12225.173
$ws.Range("I32").Value = 9792
$ws.Range("J32").Value = 18612.25
$ws.Range("K32").Value = 9792
$ws.Range("L32").Value = 18612.25
$ws.Range("M32").Value = -9505
$ws.Range("N32").Value = -19186.25

$ws.Range("H116").Value = 920850.25
$ws.Range("I116").Value = 1566.55
$ws.Range("J116").Value = 2452989.8
$ws.Range("K116").Value = 1566.55
$ws.Range("L116").Value = 2452989.8
$ws.Range("M116").Value = 727.45
$ws.Range("N116").Value = -2457577.8

$ws.Range("H122").Value = 1921.88
$ws.Range("I122").Value = 1808.6471
$ws.Range("K122").Value = 5425.9413
$ws.Range("M122").Value = -2975.9413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 920850.25
$ws.Range("I3").Value = 1566.55
$ws.Range("J3").Value = 2452989.8
$ws.Range("K3").Value = 1566.55
$ws.Range("L3").Value = 2452989.8
$ws.Range("M3").Value = -1452.55
$ws.Range("N3").Value = -2453217.8

$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("K10").Value = 2000
$ws.Range("M10").Value = -1861

$ws.Range("H59").Value = 16650.3
$ws.Range("J59").Value = 16722.555
$ws.Range("L59").Value = 16722.555
$ws.Range("N59").Value = -19012.555

$ws.Range("H99").Value = 28157.71
$ws.Range("I99").Value = 2008.9584
$ws.Range("K99").Value = 2008.9584
$ws.Range("M99").Value = -510.9584

$ws.Range("H103").Value = 9666.666999999999
$ws.Range("I103").Value = 9666.666999999999
$ws.Range("K103").Value = 9666.666999999999
$ws.Range("M103").Value = -8494.666999999999

$ws.Range("H122").Value = 1121.8334
$ws.Range("I122").Value = 761.25
$ws.Range("J122").Value = 1410.3
$ws.Range("K122").Value = 2283.75
$ws.Range("L122").Value = 4230.9
$ws.Range("M122").Value = 166.25
$ws.Range("N122").Value = -9130.9

$ws.Range("H126").Value = 28157.71
$ws.Range("I126").Value = 2008.9584
$ws.Range("K126").Value = 6026.8752
$ws.Range("M126").Value = -3556.8752

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 9372.071
$ws.Range("I120").Value = 4206.75
$ws.Range("J120").Value = 16259.167
$ws.Range("K120").Value = 12620.25
$ws.Range("L120").Value = 48777.501
$ws.Range("M120").Value = -7782.25
$ws.Range("N120").Value = -58453.501

$ws.Range("H131").Value = 4621.2144
$ws.Range("I131").Value = 7114.4443
$ws.Range("J131").Value = 4143.787
$ws.Range("K131").Value = 21343.3329
$ws.Range("L131").Value = 12431.361
$ws.Range("M131").Value = -16303.3329
$ws.Range("N131").Value = -22511.361

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 535.03705
$ws.Range("I107").Value = 368.33334
$ws.Range("J107").Value = 868.44446
$ws.Range("K107").Value = 368.33334
$ws.Range("L107").Value = 868.44446
$ws.Range("M107").Value = 1551.66666
$ws.Range("N107").Value = -4708.44446

$ws.Range("H117").Value = 8873.333000000001
$ws.Range("J117").Value = 8873.333000000001
$ws.Range("L117").Value = 8873.333000000001
$ws.Range("N117").Value = -15757.333

$ws.Range("H122").Value = 3855.158
$ws.Range("I122").Value = 5226.846
$ws.Range("J122").Value = 883.1667
$ws.Range("K122").Value = 15680.538
$ws.Range("L122").Value = 2649.5001
$ws.Range("M122").Value = -13230.538
$ws.Range("N122").Value = -7549.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1123.8823
$ws.Range("I7").Value = 1083.7273
$ws.Range("J7").Value = 1197.5
$ws.Range("K7").Value = 1083.7273
$ws.Range("L7").Value = 1197.5
$ws.Range("M7").Value = -971.7273
$ws.Range("N7").Value = -1421.5

$ws.Range("H40").Value = 2188.4211
$ws.Range("I40").Value = 1921.1538
$ws.Range("J40").Value = 2767.5
$ws.Range("K40").Value = 1921.1538
$ws.Range("L40").Value = 2767.5
$ws.Range("M40").Value = -1785.1538
$ws.Range("N40").Value = -3039.5

$ws.Range("H93").Value = 1139.9546
$ws.Range("I93").Value = 719.7368
$ws.Range("J93").Value = 3801.3333
$ws.Range("K93").Value = 719.7368
$ws.Range("L93").Value = 3801.3333
$ws.Range("M93").Value = 528.2632
$ws.Range("N93").Value = -6297.3333

$ws.Range("H118").Value = 33600
$ws.Range("J118").Value = 33600
$ws.Range("L118").Value = 33600
$ws.Range("N118").Value = -36914

$ws.Range("H126").Value = 1123.8823
$ws.Range("I126").Value = 1083.7273
$ws.Range("J126").Value = 1197.5
$ws.Range("K126").Value = 3251.1819
$ws.Range("L126").Value = 3592.5
$ws.Range("M126").Value = -781.1819
$ws.Range("N126").Value = -8532.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 321.58334
$ws.Range("I113").Value = 332.5
$ws.Range("J113").Value = 306.3
$ws.Range("K113").Value = 997.5
$ws.Range("L113").Value = 918.9000000000001
$ws.Range("M113").Value = 1172.5
$ws.Range("N113").Value = -5258.9

$ws.Range("H120").Value = 31148.75
$ws.Range("J120").Value = 31148.75
$ws.Range("L120").Value = 31148.75
$ws.Range("N120").Value = -40824.75
